$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("inputs")

# Update the scale factor in row 16 (scale_flexgrid PV) from 15 to 30
# across columns C through L.
$ws.Range("C16:L16").Value = 30
